# Auto-generated Excel COM-interop script
# Applies the cell-level value updates described in the commit diff
# for Sheets/Ultros_Profits.xlsx (workbook sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 400
$ws.Range("K12").Value = 400
$ws.Range("M12").Value = -230
$ws.Range("H69").Value = 30319.688
$ws.Range("J69").Value = 44934
$ws.Range("L69").Value = 134802
$ws.Range("N69").Value = -136550
$ws.Range("H70").Value = 5000
$ws.Range("J70").Value = 5000
$ws.Range("L70").Value = 15000
$ws.Range("N70").Value = -15540
$ws.Range("H72").Value = 30319.688
$ws.Range("J72").Value = 44934
$ws.Range("L72").Value = 404406
$ws.Range("N72").Value = -413142
$ws.Range("H73").Value = 5000
$ws.Range("J73").Value = 5000
$ws.Range("L73").Value = 15000
$ws.Range("N73").Value = -16872
$ws.Range("H74").Value = 8532.579
$ws.Range("I74").Value = 6577.625
$ws.Range("J74").Value = 9954.362999999999
$ws.Range("K74").Value = 6577.625
$ws.Range("L74").Value = 9954.362999999999
$ws.Range("M74").Value = -5641.625
$ws.Range("N74").Value = -11826.363
$ws.Range("H77").Value = 8532.579
$ws.Range("I77").Value = 6577.625
$ws.Range("J77").Value = 9954.362999999999
$ws.Range("K77").Value = 32888.125
$ws.Range("L77").Value = 49771.815
$ws.Range("M77").Value = -28208.125
$ws.Range("N77").Value = -59131.815
$ws.Range("H100").Value = 7774.3125
$ws.Range("I100").Value = 5673.625
$ws.Range("J100").Value = 9875
$ws.Range("K100").Value = 5673.625
$ws.Range("L100").Value = 9875
$ws.Range("M100").Value = -5132.625
$ws.Range("N100").Value = -10957
$ws.Range("H115").Value = 645.5
$ws.Range("I115").Value = 645.5
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1936.5
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -369.5
$ws.Range("N115").ClearContents()
$ws.Range("H132").Value = 13503.523
$ws.Range("I132").Value = 1309.6945
$ws.Range("K132").Value = 3929.0835
$ws.Range("M132").Value = -1399.0835
$ws.Range("H140").Value = 63694.25
$ws.Range("J140").Value = 95780
$ws.Range("L140").Value = 95780
$ws.Range("N140").Value = -106140
$ws.Range("H141").Value = 4993.2144
$ws.Range("I141").Value = 4790.6665
$ws.Range("K141").Value = 14371.9995
$ws.Range("M141").Value = -9191.999500000002

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 139
$ws.Range("I5").Value = 112.5
$ws.Range("J5").Value = 192
$ws.Range("K5").Value = 112.5
$ws.Range("L5").Value = 192
$ws.Range("M5").Value = -0.5
$ws.Range("N5").Value = -416
$ws.Range("H74").Value = 1197.9166
$ws.Range("I74").Value = 852.2727
$ws.Range("K74").Value = 852.2727
$ws.Range("M74").Value = 21.72730000000001
$ws.Range("H77").Value = 1197.9166
$ws.Range("I77").Value = 852.2727
$ws.Range("K77").Value = 4261.363499999999
$ws.Range("M77").Value = 106.6365000000005
$ws.Range("H135").Value = 67082.836
$ws.Range("J135").Value = 67082.836
$ws.Range("L135").Value = 67082.836
$ws.Range("N135").Value = -77222.836

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 139
$ws.Range("I4").Value = 112.5
$ws.Range("J4").Value = 192
$ws.Range("K4").Value = 112.5
$ws.Range("L4").Value = 192
$ws.Range("M4").Value = 2.5
$ws.Range("N4").Value = -422
$ws.Range("H105").Value = 3195.0908
$ws.Range("I105").Value = 3142.476
$ws.Range("K105").Value = 3142.476
$ws.Range("M105").Value = -1395.476

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 406
$ws.Range("I22").Value = 399.66666
$ws.Range("K22").Value = 399.66666
$ws.Range("M22").Value = -49.66665999999998
$ws.Range("H68").Value = 44000
$ws.Range("J68").Value = 44000
$ws.Range("L68").Value = 44000
$ws.Range("N68").Value = -45498
$ws.Range("H71").Value = 44000
$ws.Range("J71").Value = 44000
$ws.Range("L71").Value = 132000
$ws.Range("N71").Value = -139488
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H105").Value = 10421158
$ws.Range("I105").Value = 1423.7368
$ws.Range("J105").Value = 50016148
$ws.Range("K105").Value = 1423.7368
$ws.Range("L105").Value = 50016148
$ws.Range("M105").Value = 323.2632000000001
$ws.Range("N105").Value = -50019642
$ws.Range("H132").Value = 2570.08
$ws.Range("I132").Value = 2570.08
$ws.Range("K132").Value = 7710.24
$ws.Range("M132").Value = -5180.24
$ws.Range("H141").Value = 57821.395
$ws.Range("J141").Value = 57821.395
$ws.Range("L141").Value = 57821.395
$ws.Range("N141").Value = -68181.39499999999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 25001358
$ws.Range("I140").Value = 25001358
$ws.Range("K140").Value = 75004074
$ws.Range("M140").Value = -74998894

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 113201.6
$ws.Range("I80").Value = 160187.14
$ws.Range("J80").Value = 3568.6667
$ws.Range("K80").Value = 160187.14
$ws.Range("L80").Value = 3568.6667
$ws.Range("M80").Value = -159189.14
$ws.Range("N80").Value = -5564.6667
$ws.Range("H83").Value = 113201.6
$ws.Range("I83").Value = 160187.14
$ws.Range("J83").Value = 3568.6667
$ws.Range("K83").Value = 800935.7000000001
$ws.Range("L83").Value = 17843.3335
$ws.Range("M83").Value = -795943.7000000001
$ws.Range("N83").Value = -27827.3335
$ws.Range("H122").Value = 5743.5264
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 5743.5264
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 17230.5792
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -22130.5792
$ws.Range("H132").Value = 2885.0881
$ws.Range("I132").Value = 2296.0908
$ws.Range("J132").Value = 3964.9167
$ws.Range("K132").Value = 6888.2724
$ws.Range("L132").Value = 11894.7501
$ws.Range("M132").Value = -4358.2724
$ws.Range("N132").Value = -16954.7501

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 112496.75
$ws.Range("J20").Value = 112496.75
$ws.Range("L20").Value = 112496.75
$ws.Range("N20").Value = -112948.75
$ws.Range("H59").Value = 15000
$ws.Range("J59").Value = 15000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16308
$ws.Range("H68").Value = 6872.615
$ws.Range("I68").Value = 7333.3335
$ws.Range("J68").Value = 6734.4
$ws.Range("K68").Value = 7333.3335
$ws.Range("L68").Value = 6734.4
$ws.Range("M68").Value = -6584.3335
$ws.Range("N68").Value = -8232.4
$ws.Range("H71").Value = 6872.615
$ws.Range("I71").Value = 7333.3335
$ws.Range("J71").Value = 6734.4
$ws.Range("K71").Value = 36666.6675
$ws.Range("L71").Value = 33672
$ws.Range("M71").Value = -32922.6675
$ws.Range("N71").Value = -41160
$ws.Range("H102").Value = 10332.167
$ws.Range("J102").Value = 10332.167
$ws.Range("L102").Value = 10332.167
$ws.Range("N102").Value = -16822.167
$ws.Range("H136").Value = 4519.512
$ws.Range("I136").Value = 3429.0645
$ws.Range("K136").Value = 10287.1935
$ws.Range("M136").Value = -7737.193499999999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8000
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248
$ws.Range("H65").Value = 8000
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240
$ws.Range("H122").Value = 4644.421
$ws.Range("I122").Value = 3788
$ws.Range("K122").Value = 11364
$ws.Range("M122").Value = -8914
$ws.Range("H136").Value = 4128.381
$ws.Range("I136").Value = 2492.2
$ws.Range("J136").Value = 5615.8184
$ws.Range("K136").Value = 7476.599999999999
$ws.Range("L136").Value = 16847.4552
$ws.Range("M136").Value = -4926.599999999999
$ws.Range("N136").Value = -21947.4552

Write-Host "Applied Ultros_Profits.xlsx price/profit updates across all sheets."